$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 808.41174
$ws.Range("I2").Value = 237.75
$ws.Range("J2").Value = 2178
$ws.Range("K2").Value = 237.75
$ws.Range("L2").Value = 2178
$ws.Range("M2").Value = -124.75
$ws.Range("N2").Value = -2404

$ws.Range("H38").Value = 327.07144
$ws.Range("I38").Value = 125.30769
$ws.Range("J38").Value = 2950
$ws.Range("K38").Value = 375.92307
$ws.Range("L38").Value = 8850
$ws.Range("M38").Value = -3.923069999999996
$ws.Range("N38").Value = -9594

$ws.Range("H101").Value = 25001506
$ws.Range("I101").Value = 33333676
$ws.Range("K101").Value = 100001028
$ws.Range("M101").Value = -99999406

$ws.Range("H141").Value = 3999.5
$ws.Range("I141").Value = 2999.25
$ws.Range("K141").Value = 8997.75
$ws.Range("M141").Value = -3817.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2267342.2
$ws.Range("I32").Value = 2695121.2
$ws.Range("J32").Value = 877060.4
$ws.Range("K32").Value = 2695121.2
$ws.Range("L32").Value = 877060.4
$ws.Range("M32").Value = -2694834.2
$ws.Range("N32").Value = -877634.4

$ws.Range("H45").Value = 92043.09
$ws.Range("I45").Value = 167829.17
$ws.Range("J45").Value = 1099.8
$ws.Range("K45").Value = 167829.17
$ws.Range("L45").Value = 1099.8
$ws.Range("M45").Value = -167452.17
$ws.Range("N45").Value = -1853.8

$ws.Range("H61").Value = 2165.7646
$ws.Range("I61").Value = 2273.6667
$ws.Range("J61").Value = 1356.5
$ws.Range("K61").Value = 2273.6667
$ws.Range("L61").Value = 1356.5
$ws.Range("M61").Value = -2061.6667
$ws.Range("N61").Value = -1780.5

$ws.Range("H101").Value = 293301.5
$ws.Range("J101").Value = 293301.5
$ws.Range("L101").Value = 293301.5
$ws.Range("N101").Value = -299791.5

$ws.Range("H122").Value = 2891.7
$ws.Range("I122").Value = 2701.889
$ws.Range("K122").Value = 8105.667
$ws.Range("M122").Value = -5655.667

$ws.Range("H132").Value = 1638.6923
$ws.Range("I132").Value = 1564.24
$ws.Range("K132").Value = 4692.72
$ws.Range("M132").Value = -2162.72

$ws.Range("H136").Value = 2165.7646
$ws.Range("I136").Value = 2273.6667
$ws.Range("J136").Value = 1356.5
$ws.Range("K136").Value = 6821.000100000001
$ws.Range("L136").Value = 4069.5
$ws.Range("M136").Value = -4271.000100000001
$ws.Range("N136").Value = -9169.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = ""

$ws.Range("H105").Value = 2824.5
$ws.Range("I105").Value = 2849.8333
$ws.Range("J105").Value = 2748.5
$ws.Range("K105").Value = 2849.8333
$ws.Range("L105").Value = 2748.5
$ws.Range("M105").Value = -1102.8333
$ws.Range("N105").Value = -6242.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3916.3333
$ws.Range("I62").Value = 3833.3333
$ws.Range("J62").Value = 3999.3333
$ws.Range("K62").Value = 3833.3333
$ws.Range("L62").Value = 3999.3333
$ws.Range("M62").Value = -3209.3333
$ws.Range("N62").Value = -5247.3333

$ws.Range("H65").Value = 3916.3333
$ws.Range("I65").Value = 3833.3333
$ws.Range("J65").Value = 3999.3333
$ws.Range("K65").Value = 19166.6665
$ws.Range("L65").Value = 19996.6665
$ws.Range("M65").Value = -16046.6665
$ws.Range("N65").Value = -26236.6665

$ws.Range("H99").Value = 1777.6
$ws.Range("I99").Value = 1609.75
$ws.Range("J99").Value = 2449
$ws.Range("K99").Value = 1609.75
$ws.Range("L99").Value = 2449
$ws.Range("M99").Value = -111.75
$ws.Range("N99").Value = -5445

$ws.Range("H126").Value = 1777.6
$ws.Range("I126").Value = 1609.75
$ws.Range("J126").Value = 2449
$ws.Range("K126").Value = 4829.25
$ws.Range("L126").Value = 7347
$ws.Range("M126").Value = -2359.25
$ws.Range("N126").Value = -12287

$ws.Range("H132").Value = 1610.125
$ws.Range("I132").Value = 1545.7587
$ws.Range("J132").Value = 2232.3333
$ws.Range("K132").Value = 4637.2761
$ws.Range("L132").Value = 6696.999899999999
$ws.Range("M132").Value = -2107.2761
$ws.Range("N132").Value = -11756.9999

$ws.Range("H141").Value = 35991.418
$ws.Range("J141").Value = 35173.453
$ws.Range("L141").Value = 35173.453
$ws.Range("N141").Value = -45533.453

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 12127.143
$ws.Range("I70").Value = 1630
$ws.Range("K70").Value = 4890
$ws.Range("M70").Value = -4575

$ws.Range("H73").Value = 12127.143
$ws.Range("I73").Value = 1630
$ws.Range("K73").Value = 4890
$ws.Range("M73").Value = -3798

$ws.Range("H98").Value = 1392.7
$ws.Range("I98").Value = 924
$ws.Range("J98").Value = 1509.875
$ws.Range("K98").Value = 2772
$ws.Range("L98").Value = 4529.625
$ws.Range("M98").Value = -1274
$ws.Range("N98").Value = -7525.625

$ws.Range("H128").Value = 624322
$ws.Range("I128").Value = 624322
$ws.Range("K128").Value = 1872966
$ws.Range("M128").Value = -1867986

$ws.Range("H131").Value = 436704.97
$ws.Range("I131").Value = 891.7143
$ws.Range("J131").Value = 627373.25
$ws.Range("K131").Value = 2675.1429
$ws.Range("L131").Value = 1882119.75
$ws.Range("M131").Value = 2364.8571
$ws.Range("N131").Value = -1892199.75

$ws.Range("H141").Value = 6456.8
$ws.Range("I141").Value = 6456.8
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 19370.4
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -14190.4
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1500.2
$ws.Range("I122").Value = 875.25
$ws.Range("K122").Value = 2625.75
$ws.Range("M122").Value = -175.75

$ws.Range("H126").Value = 3721.75
$ws.Range("I126").Value = 3721.75
$ws.Range("K126").Value = 11165.25
$ws.Range("M126").Value = -8695.25

$ws.Range("H132").Value = 1681
$ws.Range("I132").Value = 1681
$ws.Range("K132").Value = 5043
$ws.Range("M132").Value = -2513

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6599.905
$ws.Range("J7").Value = 7293.5884
$ws.Range("L7").Value = 7293.5884
$ws.Range("N7").Value = -7517.5884

$ws.Range("H46").Value = 1794.8667
$ws.Range("I46").Value = 1092.5
$ws.Range("J46").Value = 3199.6
$ws.Range("K46").Value = 1092.5
$ws.Range("L46").Value = 3199.6
$ws.Range("M46").Value = -904.5
$ws.Range("N46").Value = -3575.6

$ws.Range("H68").Value = 2641.6667
$ws.Range("I68").Value = 2641.6667
$ws.Range("K68").Value = 2641.6667
$ws.Range("M68").Value = -1892.6667

$ws.Range("H71").Value = 2641.6667
$ws.Range("I71").Value = 2641.6667
$ws.Range("K71").Value = 13208.3335
$ws.Range("M71").Value = -9464.333500000001

$ws.Range("H126").Value = 6599.905
$ws.Range("J126").Value = 7293.5884
$ws.Range("L126").Value = 21880.7652
$ws.Range("N126").Value = -26820.7652

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19988
$ws.Range("J41").Value = 19990.334
$ws.Range("L41").Value = 19990.334
$ws.Range("N41").Value = -20770.334

$ws.Range("H126").Value = 4550
$ws.Range("I126").Value = 4550
$ws.Range("K126").Value = 13650
$ws.Range("M126").Value = -11180

$ws.Range("H132").Value = 2730.3635
$ws.Range("I132").Value = 3278.875
$ws.Range("J132").Value = 1267.6666
$ws.Range("K132").Value = 9836.625
$ws.Range("L132").Value = 3802.9998
$ws.Range("M132").Value = -7306.625
$ws.Range("N132").Value = -8862.9998
